$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.420.99"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "1.774.11"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "307.06"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "0.4269"
$ws.Range("D8").Value = "0.3618"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "0.07153"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "0.8405"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "20.43"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "1.782.63"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "6.446"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'5.250"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "0.06892"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "78.88"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "'0.000008691"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "14.89"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "26.440.45"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "5.103"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "2.011.53"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'152.20"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "1.812"
$ws.Range("E26").Value = "  -8.29%  "
$ws.Range("D27").Value = "'18.00"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "5.068"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "113.89"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "1.768"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").Value = "0.08887"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "0.7259"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "1.115"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "4.323"
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "2.746"
$ws.Range("E36").Value = "  -5.76%  "
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("D38").Value = "0.05141"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").Value = "0.1614"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "0.4912"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "2.585"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").Value = "6.342"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "'7.960"
$ws.Range("D45").Value = "104.71"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "10.13"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "1.628"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "0.06179"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "0.4454"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("D51").Value = "'1.720"
$ws.Range("E51").Value = "  +1.02%  "
